$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 39, pushing the totals block (rows 40-42) down to 41-43.
# Using Insert() on row 39 copies formatting from the row above (row 38), matching the
# styles used in the target workbook for the new row.
$ws.Rows("39:39").Insert()

# Fill in the new entry (Day 31).
$ws.Range("A39").Value = 31
$ws.Range("B39").Value = "29/5/2024"
$ws.Range("C39").Value = 3
$ws.Range("D39").Value = "Added users action logs functionality"

# Update the total formula (now on row 41) to include the new row.
$ws.Range("D41").Formula = "=SUM(C4:C39)"

# Setting a value on row 41 (which uses a 24pt font style for the "Total" label)
# causes the engine to auto-adjust the row height; restore it to avoid introducing
# an unwanted customHeight override.
$ws.Rows(41).EntireRow.AutoFit() | Out-Null

# Restore the view selection.
$ws.Range("C25").Select()
